$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.524.18"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "2.467.38"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.74%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.98%  "

$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D13").Value = "2.848.06"
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("E14").Value = "  -2.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.96%  "

$ws.Range("D16").Value = "2.464.63"
$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("E17").Value = "  -2.75%  "

$ws.Range("D18").Value = "41.517.95"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.58%  "

$ws.Range("D20").Value = "0.0$([char]0x2083)0941"
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("E29").Value = "  -1.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.50%  "

$ws.Range("E36").Value = "  -7.32%  "

$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.46%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.90%  "

$ws.Range("E41").Value = "  -5.09%  "

$ws.Range("E42").Value = "  -0.21%  "

$ws.Range("D43").Value = "1.946.01"
$ws.Range("E43").Value = "  -2.60%  "

$ws.Range("E44").Value = "  -1.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.00%  "

$ws.Range("E46").Value = "  -3.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "

$ws.Range("D48").Value = "2.706.51"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("E50").Value = "  -4.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.99%  "
